$wb = $excel.ActiveWorkbook

# "Overview" sheet: update the "Latest HO Xliff Generate Date" for the
# 7f24a935... row (row 2) to reflect a fresh handback-report generation.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-06 15:42:25"

# "zh-cn" sheet: update Correspond Handoff Datetime (H2) and
# Correspond Handback DateTime (K2) for row 2.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-06 15:42:11"
$wsZhCn.Range("K2").Value = "2016-09-06 15:42:45"

# "de-de" sheet: update Correspond Handoff Datetime (H2) and
# Correspond Handback DateTime (K2) for row 2.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-06 15:42:25"
$wsDeDe.Range("K2").Value = "2016-09-06 15:42:53"
